$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New skill entries (Growlmon / Seadramon evolutions), matching the
# existing table layout: ID | Skill | Type | Mana | CoolDown

# Row 58 -> ID 56
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = "SavageInstinct"
$ws.Range("C58").Value = "PassiveSkill(Growlmon)"
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0

# Row 59 -> ID 57
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = "ExhaustFlame"
$ws.Range("C59").Value = "DamageSkill(Growlmon)"
$ws.Range("D59").Value = 60
$ws.Range("E59").Value = 5

# Row 60 -> ID 58
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = "MagicScale"
$ws.Range("C60").Value = "PassiveSkill(Seadramon)"
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0

# Row 61 -> ID 59
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = "IceBlast"
$ws.Range("C61").Value = "DamageSkill(Seadramon)"
$ws.Range("D61").Value = 55
$ws.Range("E61").Value = 5

# Row 66 -> ID 64 (leaves rows 62-65 blank, matching the source sheet)
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = "PumpItUp"
$ws.Range("C66").Value = "StatusSkill"
$ws.Range("D66").Value = 10
$ws.Range("E66").Value = 4

# Apply the same "Good" cell style used by the rest of the data rows
$ws.Range("A58:E61").Style = "Good"
$ws.Range("A66:E66").Style = "Good"
$ws.Range("A58:E61").HorizontalAlignment = -4108
$ws.Range("A58:E61").VerticalAlignment = -4108
$ws.Range("A66:E66").HorizontalAlignment = -4108
$ws.Range("A66:E66").VerticalAlignment = -4108

# Update the view to match: scrolled down, selection on B63
$ws.Application.ActiveWindow.ScrollRow = 46
[void]$ws.Range("B63").Select()
